$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so COM does not auto-convert
# numeric-looking strings (e.g. "7.78", "1.01") into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "50.945.28"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").Value = "2.947.48"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "378.70"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D6").Value = "101.49"
$ws.Range("E6").Value = "  -1.73%  "
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "0.582"
$ws.Range("E9").Value = "  -1.55%  "
$ws.Range("D10").Value = "36.12"
$ws.Range("E10").Value = "  -1.54%  "
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").Value = "7.78"
$ws.Range("E13").Value = "  +4.71%  "
$ws.Range("D14").Value = "3.411.62"
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("D15").Value = "18.28"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("D16").Value = "12.08"
$ws.Range("E16").Value = "  +68.95%  "
$ws.Range("D17").Value = "2.948.47"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").Value = "1.01"
$ws.Range("E18").Value = "  +1.54%  "
$ws.Range("D19").Value = "50.896.56"
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("D20").Value = "3.08"
$ws.Range("E20").Value = "  -4.00%  "
$ws.Range("E21").Value = "  -1.75%  "
$ws.Range("D22").Value = "0.0₃0951"
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("D23").Value = "69.32"
$ws.Range("E23").Value = "  +1.08%  "
$ws.Range("D24").Value = "265.77"
$ws.Range("E24").Value = "  +1.19%  "
$ws.Range("E25").Value = "  +9.03%  "
$ws.Range("D26").Value = "8.08"
$ws.Range("E26").Value = "  -2.97%  "
$ws.Range("D28").Value = "7.05"
$ws.Range("E28").Value = "  -9.12%  "
$ws.Range("D29").Value = "25.59"
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("E30").Value = "  -3.66%  "
$ws.Range("E31").Value = "  -3.95%  "
$ws.Range("D32").Value = "10.16"
$ws.Range("E32").Value = "  +3.34%  "
$ws.Range("D33").Value = "50.50"
$ws.Range("E33").Value = "  -0.25%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "33.54"
$ws.Range("E35").Value = "  -2.58%  "
$ws.Range("D36").Value = "0.0432"
$ws.Range("E36").Value = "  -5.83%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "3.09"
$ws.Range("E38").Value = "  +3.08%  "
$ws.Range("E39").Value = "  +0.50%  "
$ws.Range("D40").Value = "16.55"
$ws.Range("E40").Value = "  -2.06%  "
$ws.Range("E41").Value = "  +1.25%  "
$ws.Range("D42").Value = "2.52"
$ws.Range("E42").Value = "  -2.19%  "
$ws.Range("D43").Value = "119.57"
$ws.Range("E43").Value = "  -1.46%  "
$ws.Range("E44").Value = "  +7.74%  "
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("E46").Value = "  -1.91%  "
$ws.Range("E47").Value = "  -2.22%  "
$ws.Range("D48").Value = "1.998.14"
$ws.Range("E48").Value = "  -0.97%  "
$ws.Range("E49").Value = "  -4.51%  "
$ws.Range("D50").Value = "0.0316"
$ws.Range("E50").Value = "  -10.07%  "
$ws.Range("D51").Value = "5.28"
$ws.Range("E51").Value = "  +3.68%  "
